$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.315.66"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "3.551.21"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.25"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.99"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.661"
$ws.Range("E9").Value = "  -6.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.88"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("E11").Value = "  -11.85%  "
$ws.Range("E12").Value = "  -14.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.77"
$ws.Range("E13").Value = "  -7.78%  "
$ws.Range("D14").Value = "4.151.07"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.567.53"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "66.131.39"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.01"
$ws.Range("E18").Value = "  -6.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.96"
$ws.Range("E19").Value = "  -5.37%  "
$ws.Range("E20").Value = "  -6.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.22"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.55"
$ws.Range("E23").Value = "  -5.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.82"
$ws.Range("E24").Value = "  -5.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.00"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  -5.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("E28").Value = "  -8.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("E29").Value = "  -6.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.62"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.70"
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "587.68"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.49"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.365"
$ws.Range("E39").Value = "  -7.36%  "
$ws.Range("D40").Value = "0.0₃0723"
$ws.Range("E40").Value = "  -17.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -9.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0402"
$ws.Range("E43").Value = "  -7.50%  "
$ws.Range("D44").Value = "2.730.00"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.07"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  -12.19%  "
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  -7.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.27"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.11"
$ws.Range("E50").Value = "  -9.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.52"
$ws.Range("E51").Value = "  -7.52%  "
